# Applies the "Automatic update of files" edit:
#  1. Bumps the "Förändrad" (changed) date in column C from 45171 to 45172
#     for every data row (rows 2-288).
#  2. Re-orders the two rows for case "A 36038-2023" and "A 50972-2020":
#     "A 36038-2023" now appears before "A 50972-2020" (row 5 / row 6),
#     and the "A 36038-2023" record gains an extra observed species
#     ("Plattlummer") which bumps a couple of its summary counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 288

# --- Step 1: bump column C ("Förändrad") for every data row ------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}

# --- Step 2: rewrite rows 5 and 6 with the re-ordered / updated data ----

# Row 5 becomes the "A 36038-2023" record (previously row 6), with the
# extra "Plattlummer" species added and the summary counts updated.
$ws.Range("A5").Value = "A 36038-2023"
$ws.Range("B5").Value = 45149
$ws.Range("C5").Value = 45172
$ws.Range("D5").Value = "DALARNAS LÄN"
$ws.Range("E5").Value = "ORSA"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = 4.8
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 8
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = "Smalfotad taggsvamp`r`nGarnlav`r`nMörk kolflarnlav`r`nOrange taggsvamp`r`nSvart taggsvamp`r`nTretåig hackspett`r`nVaddporing`r`nViolettgrå tagellav`r`nPlattlummer`r`nSkarp dropptaggsvamp"
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/artfynd/A 36038-2023.xlsx")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/kartor/A 36038-2023.png")'
$ws.Range("U5").Value = ""
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/klagomål/A 36038-2023.docx")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/klagomålsmail/A 36038-2023.docx")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/tillsyn/A 36038-2023.docx")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/tillsynsmail/A 36038-2023.docx")'

# Row 6 becomes the "A 50972-2020" record (previously row 5), unchanged
# apart from its new row position and the column C date bump above.
$ws.Range("A6").Value = "A 50972-2020"
$ws.Range("B6").Value = 44111
$ws.Range("C6").Value = 45172
$ws.Range("D6").Value = "DALARNAS LÄN"
$ws.Range("E6").Value = "ORSA"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 6.9
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 6
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 6
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 9
$ws.Range("R6").Value = "Kolflarnlav`r`nNordtagging`r`nVaddporing`r`nVarglav`r`nVedflamlav`r`nViolettgrå tagellav`r`nDropptaggsvamp`r`nGullgröppa`r`nSkarp dropptaggsvamp"
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/artfynd/A 50972-2020.xlsx")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/kartor/A 50972-2020.png")'
$ws.Range("U6").Value = ""
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/klagomål/A 50972-2020.docx")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/klagomålsmail/A 50972-2020.docx")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/tillsyn/A 50972-2020.docx")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ORSA/tillsynsmail/A 50972-2020.docx")'
